$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("Z1")
$rng.Borders(7).Color = 0
$rng.Borders(8).Color = 0
$rng.Borders(9).Color = 0
$rng.Font.Underline = $true
